# Updates the cryptocurrency table (rows 2-51) to match the latest scraped
# coinranking.com snapshot: refreshed prices/volumes for existing rows, and
# shifted-in a new row (BitDAO) causing every following coin to move down one
# row, with the former last row (Cronos) dropping off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as plain text (e.g. "23.132.55"), not numbers.
# Force a Text number format first so Excel does not auto-convert the strings
# we are about to assign into numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "23.132.55"
$ws.Range("E2").Value = "  -3.07%  "

# Row 3
$ws.Range("D3").Value = "1.605.24"
$ws.Range("E3").Value = "  -2.90%  "

# Row 4
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("E5").Value = "  +0.02%  "

# Row 6
$ws.Range("D6").Value = "302.41"
$ws.Range("E6").Value = "  -2.77%  "

# Row 7
$ws.Range("E7").Value = "  -2.50%  "

# Row 8
$ws.Range("D8").Value = "0.3667"
$ws.Range("E8").Value = "  -3.95%  "

# Row 9
$ws.Range("D9").Value = "50.24"
$ws.Range("E9").Value = "  -2.81%  "

# Row 10
$ws.Range("E10").Value = "  -5.06%  "

# Row 11
$ws.Range("D11").Value = "0.08161"
$ws.Range("E11").Value = "  -3.68%  "

# Row 12
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("D13").Value = "22.91"
$ws.Range("E13").Value = "  -4.47%  "

# Row 14
$ws.Range("E14").Value = "  -5.82%  "

# Row 15
$ws.Range("E15").Value = "  -3.98%  "

# Row 16
$ws.Range("D16").Value = "7.419"
$ws.Range("E16").Value = "  -7.89%  "

# Row 17
$ws.Range("D17").Value = "1.603.64"
$ws.Range("E17").Value = "  -3.05%  "

# Row 18
$ws.Range("D18").Value = "92.21"
$ws.Range("E18").Value = "  -2.08%  "

# Row 19
$ws.Range("D19").Value = "0.06880"
$ws.Range("E19").Value = "  -1.70%  "

# Row 20
$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  -6.27%  "

# Row 21
$ws.Range("D21").Value = "6.617"
$ws.Range("E21").Value = "  -5.19%  "

# Row 22
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").Value = "0.5552"
$ws.Range("E22").Value = "  -5.84%  "

# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "13.11"
$ws.Range("E24").Value = "  -4.34%  "

# Row 25
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.126.60"
$ws.Range("E25").Value = "  -3.09%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.360"
$ws.Range("E26").Value = "  -3.11%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.816"
$ws.Range("E27").Value = "  -4.80%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "21.20"
$ws.Range("E28").Value = "  -3.87%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "150.32"
$ws.Range("E29").Value = "  -2.06%  "

# Row 30
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "5.281"
$ws.Range("E30").Value = "  -2.89%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "134.14"
$ws.Range("E31").Value = "  -2.84%  "

# Row 32
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").Value = "2.390"
$ws.Range("E32").Value = "  -3.96%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.889"
$ws.Range("E33").Value = "  -11.74%  "

# Row 34
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").Value = "1.780.53"
$ws.Range("E34").Value = "  -3.21%  "

# Row 35
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.9645"
$ws.Range("E35").Value = "  -3.98%  "

# Row 36
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.07741"
$ws.Range("E36").Value = "  -4.95%  "

# Row 37
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "10.44"
$ws.Range("E37").Value = "  -3.11%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "6.311"
$ws.Range("E38").Value = "  -5.14%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02741"
$ws.Range("E39").Value = "  -5.94%  "

# Row 40
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "0.2558"
$ws.Range("E40").Value = "  -4.34%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.08908"
$ws.Range("E41").Value = "  -2.39%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.370"
$ws.Range("E42").Value = "  -4.01%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.7114"
$ws.Range("E43").Value = "  -5.79%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "12.67"
$ws.Range("E44").Value = "  -6.68%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "15.37"
$ws.Range("E45").Value = "  -6.68%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6648"
$ws.Range("E46").Value = "  -4.05%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.333"
$ws.Range("E47").Value = "  -4.66%  "

# Row 48
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "0.9993"
$ws.Range("E48").Value = "  +0.05%  "

# Row 49
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "4.008"
$ws.Range("E49").Value = "  -2.34%  "

# Row 50
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "1.258"
$ws.Range("E50").Value = "  +2.64%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "132.89"
$ws.Range("E51").Value = "  -0.16%  "
